# Backup before dimension reduction:
# Shift the "q" index labels in column A (rows 2-97) down by one,
# e.g. q1 -> q0, q2 -> q1, ..., q96 -> q95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($current -match '^q(\d+)$') {
        $n = [int]$matches[1]
        $cell.Value2 = "q$($n - 1)"
    }
}
